# Update 'Price' (D) and 'Volume(1h)' (E) columns with latest scraped values.
# Cryptos list refresh - GitHub Actions scheduled update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.240.27"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "2.625.34"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Formula = "=""597.65"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Formula = "=""152.45"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("E8").Value = "  +2.74%  "

$ws.Range("D9").Value = "2.623.49"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("E10").Value = "  +1.93%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "3.103.50"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "67.264.07"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "2.620.44"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Formula = "=""363.45"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").Formula = "=""7.49"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("E23").Value = "  +3.71%  "

$ws.Range("D24").Formula = "=""1.00"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Formula = "=""70.95"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +7.10%  "

$ws.Range("D26").Formula = "=""10.04"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").Value = "2.762.04"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").Formula = "=""1.00"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Formula = "=""575.05"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -5.75%  "

$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("D32").Formula = "=""7.82"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  -3.51%  "

$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("D38").Formula = "=""157.02"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +1.69%  "

$ws.Range("D39").Formula = "=""19.16"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").Formula = "=""5.25"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D47").Formula = "=""156.21"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)

$ws.Range("D48").Value = "0.0₆0283"
$ws.Range("E48").Value = "  -2.21%  "

$ws.Range("D49").Formula = "=""3.72"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").Formula = "=""20.55"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.90%  "

$excel.CutCopyMode = $false
